$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.190.33"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "1.585.37"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "1.807.53"
$ws.Range("E12").Value = "  -0.40%  "

$ws.Range("D13").Value = "1.588.41"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.79%  "

$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "26.177.72"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "

$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.19"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").Value = "1.413.29"
$ws.Range("E33").Value = "  +8.30%  "

$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.45%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.56%  "

$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.944"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -14.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.764"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "1.719.74"
$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.11%  "

$ws.Range("E48").Value = "  -0.94%  "

$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("E51").Value = "  -0.29%  "
